# Generate Report for Handback
# Updates row 7 ("466b1249-8b09-4fba-be63-7dcbd48901de") on the zh-cn and
# de-de sheets: a handback file was produced, but it was not based on the
# latest handoff, so the target file / handback datetime / error detail
# columns now get populated (instead of being blank / the default
# 0001-01-01 placeholder).

$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6cb408a7ac9ac7381966f160067eae2cf996671c/e2e/466b1249-8b09-4fba-be63-7dcbd48901de.md"
$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7f5ffd76f2f1911778ebbddc7c2919ace8b537b/e2e/466b1249-8b09-4fba-be63-7dcbd48901de.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6cb408a7ac9ac7381966f160067eae2cf996671c/e2e/466b1249-8b09-4fba-be63-7dcbd48901de.md."

function Update-HandbackRow7 {
    param(
        $ws,
        [string]$targetFileValue,
        [string]$handbackDateTime
    )

    # I7: Latest Target File - now has a value and becomes a hyperlink
    $ws.Range("I7").Value = "466b1249-8b09-4fba-be63-7dcbd48901de.md"
    $ws.Hyperlinks.Add($ws.Range("I7"), $latestUrl, "", "", "466b1249-8b09-4fba-be63-7dcbd48901de.md")
    $ws.Range("I7").Font.Underline = $true
    $ws.Range("I7").Font.Color = 15570276

    # J7: Latest Handback File
    $ws.Range("J7").Value = $targetFileValue

    # K7: Latest Handback DateTime
    $ws.Range("K7").Value = $handbackDateTime

    # P7: Error Detail
    $ws.Range("P7").Value = $errorMessage
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow7 -ws $wsZhCn `
    -targetFileValue "466b1249-8b09-4fba-be63-7dcbd48901de.88cc4bceca17b5a925011e4b4c6c8a5f6491fa78.zh-cn.xlf" `
    -handbackDateTime "2016-08-30 06:58:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow7 -ws $wsDeDe `
    -targetFileValue "466b1249-8b09-4fba-be63-7dcbd48901de.88cc4bceca17b5a925011e4b4c6c8a5f6491fa78.de-de.xlf" `
    -handbackDateTime "2016-08-30 06:58:35"
